# Refresh market-price-derived columns (H..N: currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) on the rows
# touched by the scheduled market-data runner, across all eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 97.833336
$ws.Range("I5").Value = 88
$ws.Range("J5").Value = 117.5
$ws.Range("K5").Value = 88
$ws.Range("L5").Value = 117.5
$ws.Range("M5").Value = 27
$ws.Range("N5").Value = -347.5
# Row 10
$ws.Range("H10").Value = 8399
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 8399
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 8399
$ws.Range("N10").Value = -8985
# Row 11
$ws.Range("H11").Value = 87.53846
$ws.Range("I11").Value = 87.53846
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 87.53846
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 52.46154
# Row 17
$ws.Range("H17").Value = 929.6
$ws.Range("I17").Value = 864.0526
$ws.Range("J17").Value = 1007.4375
$ws.Range("K17").Value = 2592.1578
$ws.Range("L17").Value = 3022.3125
$ws.Range("M17").Value = -2424.1578
$ws.Range("N17").Value = -3358.3125
# Row 33
$ws.Range("H33").Value = 209.6875
$ws.Range("I33").Value = 144.33333
$ws.Range("J33").Value = 405.75
$ws.Range("K33").Value = 144.33333
$ws.Range("L33").Value = 405.75
$ws.Range("M33").Value = 84.66667000000001
$ws.Range("N33").Value = -863.75
# Row 51
$ws.Range("H51").Value = 20000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 20000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -20968
# Row 113
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1746
$ws.Range("N113").Value = ""
# Row 132
$ws.Range("H132").Value = 1449.1428
$ws.Range("I132").Value = 1789.6
$ws.Range("J132").Value = 598
$ws.Range("K132").Value = 5368.799999999999
$ws.Range("L132").Value = 1794
$ws.Range("M132").Value = -2838.799999999999
$ws.Range("N132").Value = -6854
# Row 137
$ws.Range("H137").Value = 2008.7858
$ws.Range("I137").Value = 1943.125
$ws.Range("J137").Value = 2096.3333
$ws.Range("K137").Value = 5829.375
$ws.Range("L137").Value = 6288.999899999999
$ws.Range("M137").Value = -3279.375
$ws.Range("N137").Value = -11388.9999
# Row 138
$ws.Range("H138").Value = 2311.463
$ws.Range("I138").Value = 1901.56
$ws.Range("J138").Value = 2664.8276
$ws.Range("K138").Value = 5704.68
$ws.Range("L138").Value = 7994.4828
$ws.Range("M138").Value = -564.6800000000003
$ws.Range("N138").Value = -18274.4828

$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 18499.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 18499.5
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 18499.5
$ws.Range("N23").Value = -19017.5
# Row 37
$ws.Range("H37").Value = 14611.111
$ws.Range("I37").Value = 5000
$ws.Range("J37").Value = 22300
$ws.Range("K37").Value = 5000
$ws.Range("L37").Value = 22300
$ws.Range("M37").Value = -4727
$ws.Range("N37").Value = -22846
# Row 55
$ws.Range("H55").Value = 21991.166
$ws.Range("I55").Value = 12024
$ws.Range("J55").Value = 26974.75
$ws.Range("K55").Value = 12024
$ws.Range("L55").Value = 26974.75
$ws.Range("M55").Value = -11709
$ws.Range("N55").Value = -27604.75
# Row 61
$ws.Range("H61").Value = 1937.08
$ws.Range("I61").Value = 2028.5454
$ws.Range("J61").Value = 1266.3334
$ws.Range("K61").Value = 2028.5454
$ws.Range("L61").Value = 1266.3334
$ws.Range("M61").Value = -1816.5454
$ws.Range("N61").Value = -1690.3334
# Row 122
$ws.Range("H122").Value = 7019.7
$ws.Range("I122").Value = 1585.4286
$ws.Range("J122").Value = 19699.666
$ws.Range("K122").Value = 4756.2858
$ws.Range("L122").Value = 59098.99800000001
$ws.Range("M122").Value = -2306.2858
$ws.Range("N122").Value = -63998.99800000001
# Row 132
$ws.Range("H132").Value = 2999.5
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -17057
# Row 136
$ws.Range("H136").Value = 1937.08
$ws.Range("I136").Value = 2028.5454
$ws.Range("J136").Value = 1266.3334
$ws.Range("K136").Value = 6085.6362
$ws.Range("L136").Value = 3799.0002
$ws.Range("M136").Value = -3535.6362
$ws.Range("N136").Value = -8899.0002

$ws = $wb.Worksheets.Item("BSM")
# Row 33
$ws.Range("H33").Value = 11355.571
$ws.Range("I33").Value = 2994.5
$ws.Range("J33").Value = 14700
$ws.Range("K33").Value = 2994.5
$ws.Range("L33").Value = 14700
$ws.Range("M33").Value = -2658.5
$ws.Range("N33").Value = -15372
# Row 87
$ws.Range("H87").Value = 51833.332
$ws.Range("I87").Value = 75000
$ws.Range("J87").Value = 47200
$ws.Range("K87").Value = 75000
$ws.Range("L87").Value = 47200
$ws.Range("M87").Value = -73752
$ws.Range("N87").Value = -49696
# Row 90
$ws.Range("H90").Value = 51833.332
$ws.Range("I90").Value = 75000
$ws.Range("J90").Value = 47200
$ws.Range("K90").Value = 225000
$ws.Range("L90").Value = 141600
$ws.Range("M90").Value = -218760
$ws.Range("N90").Value = -154080
# Row 107
$ws.Range("H107").Value = 1807.2727
$ws.Range("I107").Value = 1775.7778
$ws.Range("J107").Value = 1949
$ws.Range("K107").Value = 1775.7778
$ws.Range("L107").Value = 1949
$ws.Range("M107").Value = 144.2221999999999
$ws.Range("N107").Value = -5789
# Row 133
$ws.Range("H133").Value = 99944
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 99944
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 99944
$ws.Range("N133").Value = -110064
# Row 134
$ws.Range("H134").Value = 2084.9167
$ws.Range("I134").Value = 1502
$ws.Range("J134").Value = 4999.5
$ws.Range("K134").Value = 4506
$ws.Range("L134").Value = 14998.5
$ws.Range("M134").Value = -1971
$ws.Range("N134").Value = -20068.5

$ws = $wb.Worksheets.Item("CRP")
# Row 32
$ws.Range("H32").Value = 351
$ws.Range("I32").Value = 351
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 351
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -35
# Row 59
$ws.Range("H59").Value = 49331
$ws.Range("I59").Value = 37995
$ws.Range("J59").Value = 54999
$ws.Range("K59").Value = 37995
$ws.Range("L59").Value = 54999
$ws.Range("M59").Value = -36850
$ws.Range("N59").Value = -57289
# Row 132
$ws.Range("H132").Value = 1596.96
$ws.Range("I132").Value = 1173.0588
$ws.Range("J132").Value = 2497.75
$ws.Range("K132").Value = 3519.1764
$ws.Range("L132").Value = 7493.25
$ws.Range("M132").Value = -989.1764000000003
$ws.Range("N132").Value = -12553.25

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 192.64285
$ws.Range("I12").Value = 158.7
$ws.Range("J12").Value = 277.5
$ws.Range("K12").Value = 476.1
$ws.Range("L12").Value = 832.5
$ws.Range("M12").Value = -303.1
$ws.Range("N12").Value = -1178.5
# Row 51
$ws.Range("H51").Value = 1317.375
$ws.Range("I51").Value = 1310.8
$ws.Range("J51").Value = 1328.3334
$ws.Range("K51").Value = 3932.4
$ws.Range("L51").Value = 3985.0002
$ws.Range("M51").Value = -3472.4
$ws.Range("N51").Value = -4905.0002
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = ""

$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 10944.777
$ws.Range("I7").Value = 7505
$ws.Range("J7").Value = 11927.571
$ws.Range("K7").Value = 7505
$ws.Range("L7").Value = 11927.571
$ws.Range("M7").Value = -7393
$ws.Range("N7").Value = -12151.571
# Row 8
$ws.Range("H8").Value = 10944.777
$ws.Range("I8").Value = 7505
$ws.Range("J8").Value = 11927.571
$ws.Range("K8").Value = 7505
$ws.Range("L8").Value = 11927.571
$ws.Range("M8").Value = -7366
$ws.Range("N8").Value = -12205.571
# Row 11
$ws.Range("H11").Value = 6669999.5
$ws.Range("I11").Value = 6669999.5
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 6669999.5
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -6669860.5
$ws.Range("N11").Value = ""
# Row 80
$ws.Range("H80").Value = 5794.125
$ws.Range("I80").Value = 5283.3335
$ws.Range("J80").Value = 6100.6
$ws.Range("K80").Value = 5283.3335
$ws.Range("L80").Value = 6100.6
$ws.Range("M80").Value = -4285.3335
$ws.Range("N80").Value = -8096.6
# Row 83
$ws.Range("H83").Value = 5794.125
$ws.Range("I83").Value = 5283.3335
$ws.Range("J83").Value = 6100.6
$ws.Range("K83").Value = 26416.6675
$ws.Range("L83").Value = 30503
$ws.Range("M83").Value = -21424.6675
$ws.Range("N83").Value = -40487
# Row 122
$ws.Range("H122").Value = 68987.664
$ws.Range("I122").Value = 1801.7273
$ws.Range("J122").Value = 253749
$ws.Range("K122").Value = 5405.1819
$ws.Range("L122").Value = 761247
$ws.Range("M122").Value = -2955.1819
$ws.Range("N122").Value = -766147
# Row 126
$ws.Range("H126").Value = 4522.273
$ws.Range("I126").Value = 4335.2856
$ws.Range("J126").Value = 4849.5
$ws.Range("K126").Value = 13005.8568
$ws.Range("L126").Value = 14548.5
$ws.Range("M126").Value = -10535.8568
$ws.Range("N126").Value = -19488.5
# Row 132
$ws.Range("H132").Value = 940.7273
$ws.Range("I132").Value = 634.8
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 1904.4
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = 625.6000000000001
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 19152.5
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 19152.5
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 19152.5
$ws.Range("N14").Value = -19496.5
# Row 26
$ws.Range("H26").Value = 17499.5
$ws.Range("I26").Value = 14999
$ws.Range("J26").Value = 20000
$ws.Range("K26").Value = 14999
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = -14704
$ws.Range("N26").Value = -20590
# Row 32
$ws.Range("H32").Value = 10001
$ws.Range("I32").Value = 10001
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 10001
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9684
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = ""
# Row 43
$ws.Range("H43").Value = 9614.857
$ws.Range("I43").Value = 6505.5
$ws.Range("J43").Value = 10858.6
$ws.Range("K43").Value = 6505.5
$ws.Range("L43").Value = 10858.6
$ws.Range("M43").Value = -6312.5
$ws.Range("N43").Value = -11244.6
# Row 82
$ws.Range("H82").Value = 900.58826
$ws.Range("I82").Value = 608.6923
$ws.Range("J82").Value = 1849.25
$ws.Range("K82").Value = 608.6923
$ws.Range("L82").Value = 1849.25
$ws.Range("M82").Value = -247.6923
$ws.Range("N82").Value = -2571.25
# Row 85
$ws.Range("H85").Value = 900.58826
$ws.Range("I85").Value = 608.6923
$ws.Range("J85").Value = 1849.25
$ws.Range("K85").Value = 608.6923
$ws.Range("L85").Value = 1849.25
$ws.Range("M85").Value = 639.3077
$ws.Range("N85").Value = -4345.25
# Row 93
$ws.Range("H93").Value = 892.6087
$ws.Range("I93").Value = 750.6923
$ws.Range("J93").Value = 1077.1
$ws.Range("K93").Value = 750.6923
$ws.Range("L93").Value = 1077.1
$ws.Range("M93").Value = 497.3077
$ws.Range("N93").Value = -3573.1
# Row 122
$ws.Range("H122").Value = 3748.3333
$ws.Range("I122").Value = 3747.5
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 11242.5
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -8792.5
$ws.Range("N122").Value = -16150
# Row 132
$ws.Range("H132").Value = 4578.8
$ws.Range("I132").Value = 4578.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13736.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11206.4
$ws.Range("N132").Value = ""
# Row 136
$ws.Range("H136").Value = 6328.5884
$ws.Range("I136").Value = 5907.0713
$ws.Range("J136").Value = 8295.666999999999
$ws.Range("K136").Value = 17721.2139
$ws.Range("L136").Value = 24887.001
$ws.Range("M136").Value = -15171.2139
$ws.Range("N136").Value = -29987.001

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 4995
$ws.Range("I126").Value = 4995
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14985
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12515
# Row 132
$ws.Range("H132").Value = 2733.5334
$ws.Range("I132").Value = 1958.6666
$ws.Range("J132").Value = 5833
$ws.Range("K132").Value = 5875.9998
$ws.Range("L132").Value = 17499
$ws.Range("M132").Value = -3345.9998
$ws.Range("N132").Value = -22559
# Row 136
$ws.Range("H136").Value = 1301.381
$ws.Range("I136").Value = 851.6111
$ws.Range("J136").Value = 2554.8333
$ws.Range("K136").Value = 2554.8333
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -4.833299999999781
$ws.Range("N136").Value = -17100
